# Updated cryptos list on Fri Mar  8 07:58:24 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.246.34"
$ws.Range("E2").Value = "  +1.44%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.915.59"
$ws.Range("E3").Value = "  +3.34%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "468.06"
$ws.Range("E5").Value = "  +8.12%  "

# Row 6 - Solana
$ws.Range("D6").Value = "144.71"
$ws.Range("E6").Value = "  +5.25%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  -0.71%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.05%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.70%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +8.74%  "

# Row 11 - ShibaInu
$ws.Range("E11").Value = "  +8.86%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "43.17"
$ws.Range("E12").Value = "  +1.34%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.550.60"
$ws.Range("E13").Value = "  +3.74%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "10.35"
$ws.Range("E14").Value = "  -0.47%  "

# Row 15 - Uniswap
$ws.Range("D15").Value = "14.92"
$ws.Range("E15").Value = "  -0.13%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.941.51"
$ws.Range("E16").Value = "  +3.98%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  -0.33%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "19.82"
$ws.Range("E18").Value = "  -0.69%  "

# Row 19 - Polygon
$ws.Range("D19").Value = "1.15"
$ws.Range("E19").Value = "  +2.41%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "67.508.22"
$ws.Range("E20").Value = "  +1.71%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "431.03"
$ws.Range("E21").Value = "  +5.93%  "

# Row 22 - InternetComputer(DFINITY)
$ws.Range("D22").Value = "14.64"
$ws.Range("E22").Value = "  -0.83%  "

# Row 23 - ImmutableX
$ws.Range("D23").Value = "3.35"
$ws.Range("E23").Value = "  +2.86%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "87.49"
$ws.Range("E24").Value = "  +3.08%  "

# Row 25 - was EthereumClassic, now PancakeSwap
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "3.55"
$ws.Range("E25").Value = "  +6.62%  "

# Row 26 - was PancakeSwap, now EthereumClassic
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "38.41"
$ws.Range("E26").Value = "  +4.01%  "

# Row 27 - Filecoin
$ws.Range("D27").Value = "10.29"
$ws.Range("E27").Value = "  +4.28%  "

# Row 28 - LEO
$ws.Range("D28").Value = "5.74"
$ws.Range("E28").Value = "  +3.58%  "

# Row 29 - RenderToken
$ws.Range("D29").Value = "9.63"
$ws.Range("E29").Value = "  -0.96%  "

# Row 30 - Bittensor
$ws.Range("D30").Value = "727.18"
$ws.Range("E30").Value = "  +3.16%  "

# Row 31 - Cosmos
$ws.Range("D31").Value = "13.53"
$ws.Range("E31").Value = "  -1.88%  "

# Row 32 - Hedera
$ws.Range("E32").Value = "  -4.21%  "

# Row 33 - Toncoin
$ws.Range("E33").Value = "  +4.78%  "

# Row 34 - InjectiveProtocol
$ws.Range("D34").Value = "42.95"
$ws.Range("E34").Value = "  +3.29%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  +4.19%  "

# Row 36 - OKB
$ws.Range("D36").Value = "57.84"
$ws.Range("E36").Value = "  +2.95%  "

# Row 37 - Dai
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.15%  "

# Row 38 - PEPE
$ws.Range("D38").Value = "0.0" + [char]0x2083 + "0790"
$ws.Range("E38").Value = "  +17.58%  "

# Row 39 - NEARProtocol
$ws.Range("D39").Value = "5.36"
$ws.Range("E39").Value = "  -2.87%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +1.24%  "

# Row 41 - ThetaToken
$ws.Range("D41").Value = "3.03"
$ws.Range("E41").Value = "  +3.71%  "

# Row 42 - Fetch.AI
$ws.Range("D42").Value = "2.63"
$ws.Range("E42").Value = "  -4.39%  "

# Row 43 - Stellar
$ws.Range("E43").Value = "  -0.41%  "

# Row 44 - TheGraph
$ws.Range("D44").Value = "0.338"
$ws.Range("E44").Value = "  +3.07%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  +0.07%  "

# Row 46 - WEMIXToken
$ws.Range("E46").Value = "  +4.22%  "

# Row 47 - ARBITRUM
$ws.Range("D47").Value = "2.17"
$ws.Range("E47").Value = "  +4.84%  "

# Row 48 - LidoDAOToken
$ws.Range("D48").Value = "3.39"
$ws.Range("E48").Value = "  +1.52%  "

# Row 49 - Monero
$ws.Range("D49").Value = "146.17"
$ws.Range("E49").Value = "  +4.13%  "

# Row 50 - ApeXProtocol
$ws.Range("D50").Value = "3.13"
$ws.Range("E50").Value = "  -2.87%  "

# Row 51 - Stacks
$ws.Range("E51").Value = "  +3.11%  "
